$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.487.69"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").Value = "2.458.97"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.51%  "
$ws.Range("D9").Value = "2.457.84"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("E10").Value = "  -5.48%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.84%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "2.904.82"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "68.322.75"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.00%  "
$ws.Range("D18").Value = "2.438.61"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("E25").Value = "  -4.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.74%  "
$ws.Range("E27").Value = "  -4.87%  "
$ws.Range("D28").Value = "2.583.97"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.58%  "
$ws.Range("D30").Value = "0.0₃0839"
$ws.Range("E30").Value = "  -5.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +136.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "435.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  -4.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("E48").Value = "  -3.24%  "
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("E50").Value = "  -6.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.563"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.19%  "

Write-Output "edits applied"
